$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.482.65"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.933.78"
$ws.Range("E3").Value = "  +4.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.27"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4739"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2868"
$ws.Range("E8").Value = "  +4.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06642"
$ws.Range("E9").Value = "  +4.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "107.90"
$ws.Range("E10").Value = "  +27.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.90"
$ws.Range("E11").Value = "  +5.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.925.30"
$ws.Range("E12").Value = "  +3.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07624"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.168"
$ws.Range("E14").Value = "  +3.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6603"
$ws.Range("E15").Value = "  +5.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "305.55"
$ws.Range("E16").Value = "  +20.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.507.01"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.02"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007571"
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.177.50"
$ws.Range("E21").Value = "  +3.91%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.271"
$ws.Range("E23").Value = "  +6.86%  "
$ws.Range("E24").Value = "  +6.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.44"
$ws.Range("E25").Value = "  +2.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.292"
$ws.Range("E26").Value = "  +3.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.84"
$ws.Range("E27").Value = "  +15.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.052"
$ws.Range("E28").Value = "  +9.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1107"
$ws.Range("E29").Value = "  +8.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.361"
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.089"
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.947"
$ws.Range("E32").Value = "  +2.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05025"
$ws.Range("E33").Value = "  +3.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7417"
$ws.Range("E34").Value = "  +6.08%  "
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.756"
$ws.Range("E36").Value = "  +2.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01969"
$ws.Range("E37").Value = "  +4.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.689"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.037"
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8828"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.68"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("E42").Value = "  +11.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.805"
$ws.Range("E43").Value = "  +5.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4182"
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.251"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.224"
$ws.Range("E47").Value = "  +7.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1214"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.87"
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05630"
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3848"
$ws.Range("E51").Value = "  +4.22%  "
